$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 2647
$ws1.Range("F10").Value = 1631
$ws1.Range("F11").Value = 1661
$ws1.Range("F14").Value = 680
$ws1.Range("F15").Value = 847
$ws1.Range("F16").Value = 124
$ws1.Range("F17").Value = 347
$ws1.Range("F18").Value = 1102
$ws1.Range("F22").Value = 5863
$ws1.Range("F23").Value = 239
$ws1.Range("F24").Value = 1106
$ws1.Range("F27").Value = 149
$ws1.Range("F28").Value = 271
$ws1.Range("F30").Value = 50
$ws1.Range("F31").Value = 1069
$ws1.Range("F32").Value = 849
$ws1.Range("F34").Value = 73
$ws1.Range("F37").Value = 1224
$ws1.Range("F39").Value = 125
$ws1.Range("F42").Value = 138

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F14").Value = 2647
$ws4.Range("F15").Value = 1631
$ws4.Range("F16").Value = 1661
$ws4.Range("F19").Value = 680
$ws4.Range("F21").Value = 847
$ws4.Range("F22").Value = 124
$ws4.Range("F23").Value = 347
$ws4.Range("F24").Value = 1102
$ws4.Range("F27").Value = 5863
$ws4.Range("F28").Value = 239
$ws4.Range("F29").Value = 1106
$ws4.Range("F32").Value = 149
$ws4.Range("F33").Value = 271
$ws4.Range("F35").Value = 50
$ws4.Range("F36").Value = 1069
$ws4.Range("F37").Value = 849
$ws4.Range("F39").Value = 73
$ws4.Range("F41").Value = 1224
$ws4.Range("F43").Value = 125
$ws4.Range("F46").Value = 138

$wb.Save()
